# Case_5_14 res_bus/vm_pu.xlsx update: case with 380 kV done.
# Bus 0/1 setpoint voltage (column B) moves from 1.05 pu to 1.02 pu, and the
# resulting bus voltage magnitudes (columns C:F and I:M) are refreshed with
# the newly computed power-flow results for rows 2-25 (buses 0-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B:F (per data row, top to bottom)
$rowsBF = @(
    @(1.02,1.080686369548288,1.081314231955256,1.083625367113356,1.093730810571756),
    @(1.02,1.081900152883504,1.082279443728756,1.084690541199117,1.094823840609915),
    @(1.02,1.082685458684608,1.082903842965749,1.085379832326296,1.095531235005786),
    @(1.02,1.083015580419484,1.083166303299428,1.085669623548777,1.095828655638714),
    @(1.02,1.083071008133466,1.083210369384687,1.08571828156711,1.095878595770214),
    @(1.02,1.082689869867341,1.082907350118671,1.085383704477892,1.095535209028109),
    @(1.02,1.081096593495043,1.081640462703724,1.083985337293444,1.094100178512696),
    @(1.02,1.078288256736938,1.079406814918399,1.081521603372533,1.091572442080405),
    @(1.02,1.076415412255023,1.077916851078964,1.079879313134068,1.089887890111294),
    @(1.02,1.075604280900414,1.077271465728614,1.079168218982717,1.089158592844756),
    @(1.02,1.075302961965501,1.077031706740819,1.078904090117165,1.088887717047011),
    @(1.02,1.075367597227633,1.0770831374167,1.078960746529895,1.088945820042463),
    @(1.02,1.075579374378904,1.077251647877423,1.079146385950103,1.089136201812411),
    @(1.02,1.075709853417458,1.077355468183827,1.079260764910172,1.089253504614627),
    @(1.02,1.076469240465385,1.077959678442116,1.079926506614846,1.089936293707444),
    @(1.02,1.076945535014761,1.07833862372198,1.080344115616414,1.09036462193219),
    @(1.02,1.077223332590492,1.078559634812967,1.08058770269156,1.090614470657842),
    @(1.02,1.077318051565926,1.078634990301643,1.08067076009929,1.09069966470498),
    @(1.02,1.076894434887316,1.078297968687141,1.080299309839842,1.090318665121411),
    @(1.02,1.0755170120927,1.077202026696253,1.079091719689272,1.089080138646943),
    @(1.02,1.074650805057516,1.076512766537886,1.078332477593226,1.088301529645892),
    @(1.02,1.075110014258162,1.076878175452934,1.078734964876799,1.088714275724995),
    @(1.02,1.076917524887247,1.078316339012849,1.080319555632598,1.090339430988218),
    @(1.02,1.079014381391029,1.079984416128499,1.082158499029559,1.092225812340163)
)

# New values for columns I:M (per data row, top to bottom)
$rowsIM = @(
    @(1.067001903477116,1.085564674350209,1.083985378331075,1.086290496804683,1.096369954386103),
    @(1.06744291583022,1.086438655726869,1.084768655942776,1.087173930016216,1.097283056065761),
    @(1.067726974172488,1.087003501696314,1.085274706203559,1.087745016320181,1.097873421732351),
    @(1.067846080079299,1.087240801340488,1.085487262998474,1.087984969089828,1.098121499118141),
    @(1.067866060190167,1.087280635542253,1.085522941288429,1.08802525050155,1.098163145835497),
    @(1.067728566897663,1.087006673139733,1.085277547130776,1.087748223098242,1.097876736995062),
    @(1.067151216395424,1.085860181410225,1.084250253129278,1.086589172355008,1.096678639562403),
    @(1.066123826274001,1.083834683858763,1.082434005563909,1.084542500478427,1.09456379083628),
    @(1.065432128266147,1.082480780453503,1.081219080078009,1.083175129020921,1.093151392348076),
    @(1.065131001947355,1.081893665542158,1.080692022600122,1.082582337180067,1.092539204261199),
    @(1.065018906653415,1.081675453662155,1.080496100757731,1.082362039995718,1.092311717611416),
    @(1.065042962507966,1.081722266809985,1.08053813342667,1.082409299389723,1.09236051848354),
    @(1.065121741088266,1.081875630764609,1.080675830689405,1.082564129553498,1.092520402041682),
    @(1.065170246874864,1.081970105929991,1.080760650713998,1.082659511197131,1.092618899258527),
    @(1.06545207890957,1.082519726994063,1.081254038273877,1.083214455591063,1.093192008308268),
    @(1.065628431411616,1.082864257131753,1.081563262386787,1.083562366434179,1.093551340354096),
    @(1.065731138919414,1.083065132274109,1.081743532441508,1.083765228416611,1.093760873850521),
    @(1.065766133099684,1.083133611386652,1.081804983764316,1.083834387519954,1.093832309415538),
    @(1.065609526589977,1.082827300949389,1.081530095413515,1.083525045988942,1.093512793516351),
    @(1.065098549489319,1.08183047253541,1.080635286404863,1.082518538916599,1.092473322887284),
    @(1.064775868326183,1.081202966927916,1.080071820592793,1.081885083155689,1.091819229376715),
    @(1.064947061564191,1.081535691959733,1.080370606694718,1.082220949627523,1.092166027983998),
    @(1.065618069346679,1.082844000114949,1.081545082434729,1.083541909703098,1.093530211354463),
    @(1.066390622446641,1.084358949137674,1.082904267470118,1.085072125968345,1.095110967029838)
)

function ConvertTo-Excel2DArray($rows) {
    $nRows = $rows.Count
    $nCols = $rows[0].Count
    $arr = New-Object 'object[,]' $nRows, $nCols
    for ($i = 0; $i -lt $nRows; $i++) {
        for ($j = 0; $j -lt $nCols; $j++) {
            $arr[$i, $j] = $rows[$i][$j]
        }
    }
    return $arr
}

$ws.Range("B2:F25").Value = (ConvertTo-Excel2DArray $rowsBF)
$ws.Range("I2:M25").Value = (ConvertTo-Excel2DArray $rowsIM)
